$wb = $excel.ActiveWorkbook
try { Write-Output "AppWidth=$($excel.Width)" } catch { Write-Output "ERR $_" }
try { Write-Output "AppLeft=$($excel.Left)" } catch { Write-Output "ERR $_" }
$win = $wb.Windows.Item(1)
try { Write-Output "WinWidth=$($win.Width)" } catch { Write-Output "ERR $_" }
try { Write-Output "WinLeft=$($win.Left)" } catch { Write-Output "ERR $_" }
try { Write-Output "WinTop=$($win.Top)" } catch { Write-Output "ERR $_" }
